$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------------
# Sheet "Smoke_Suite" (Worksheets(1)) - main test-data grid
# ---------------------------------------------------------------------------
$ws1 = $wb.Worksheets.Item(1)

# Header row: "Imact" was a typo for "Impact"
$ws1.Range("I1").Value = "Impact"

# Row 2 - rename test method to the lower-cased java convention, bump the
# Problem ticket number, and correct the Impact/Complexity picklist values
$ws1.Range("B2").Value = "testCreateProblemTicket"
$ws1.Range("C2").Value = "PRB02001662"
$ws1.Range("I2").Value = "3 - Moderate / Limited"
$ws1.Range("J2").Value = "3 - Low"

# Row 3 - rename test method, fill in Problem ID + Status that were blank,
# and correct the Impact/Complexity picklist values
$ws1.Range("B3").Value = "testUpdateProblemTicket"
$ws1.Range("C3").Value = "PRB02001662"
$ws1.Range("E3").Value = "Passed"
$ws1.Range("I3").Value = "2 - Significant / Large"
$ws1.Range("J3").Value = "2 - Medium"
$ws1.Range("C3").Style = "Normal"
$ws1.Range("E3").Style = "Normal"

# Row 4 - test method replaced (approve -> "different phases"), Problem ID,
# Phase and Status filled in
$ws1.Range("B4").Value = "testDifferentPhasesOfProblemTicket"
$ws1.Range("C4").Value = "PRB02001662"
$ws1.Range("D4").Value = "Accepted"
$ws1.Range("E4").Value = "Passed"
$ws1.Range("C4:E4").Style = "Normal"

# ---------------------------------------------------------------------------
# Selection / view tweaks
# ---------------------------------------------------------------------------
$ws1.Activate()
$ws1.Range("C3").Select()
